# 2022 January 8 - Minor fixes in modeling code.
# Correct the numeric placeholders in the cost-function confusion-matrix
# table on "Feuil1": the TN / FP cost values used "2898" where "28983"
# was intended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("K9").Value  = "TN *( +28983)"
$ws.Range("K10").Value = " FP * (-28983)"

# Restore the active selection to match the author's session.
$ws.Range("L20").Select() | Out-Null
